$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table cell-by-cell to match the latest scrape.
# Values that Excel would otherwise auto-parse as numbers/dates are written with a
# leading apostrophe (quote-prefix) so they stay literal text, exactly as scraped
# (this preserves things like trailing zeros, e.g. '68.70' instead of 68.7).

$ws.Range("D2").Value = "27.530.77"
$ws.Range("E2").Value = "  +4.97%  "
$ws.Range("D3").Value = "1.723.14"
$ws.Range("E3").Value = "  +4.05%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'229.67"
$ws.Range("E5").Value = "  +4.72%  "
$ws.Range("D6").Value = "'0.5413"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.2757"
$ws.Range("D9").Value = "'0.06827"
$ws.Range("E9").Value = "  +7.29%  "
$ws.Range("D10").Value = "'21.56"
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("D11").Value = "'0.07778"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "'4.723"
$ws.Range("E12").Value = "  +3.10%  "
$ws.Range("D13").Value = "1.753.41"
$ws.Range("E13").Value = "  +6.00%  "
$ws.Range("D14").Value = "1.958.97"
$ws.Range("E14").Value = "  +3.87%  "
$ws.Range("D15").Value = "'0.5985"
$ws.Range("E15").Value = "  +5.54%  "
$ws.Range("D16").Value = "0.0₅8412"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").Value = "'68.70"
$ws.Range("E17").Value = "  +4.86%  "
$ws.Range("D18").Value = "27.473.09"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").Value = "'4.807"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'210.36"
$ws.Range("E21").Value = "  +9.41%  "
$ws.Range("D22").Value = "'10.92"
$ws.Range("E22").Value = "  +5.02%  "
$ws.Range("D23").Value = "'6.232"
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'146.29"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").Value = "'0.1252"
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("D27").Value = "'7.443"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").Value = "'16.89"
$ws.Range("E28").Value = "  +5.28%  "
$ws.Range("D29").Value = "'1.622"
$ws.Range("E29").Value = "  +8.78%  "
$ws.Range("D30").Value = "'0.05593"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("D32").Value = "'3.670"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").Value = "'3.523"
$ws.Range("E33").Value = "  +4.40%  "
$ws.Range("D34").Value = "'1.629"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").Value = "'0.9771"
$ws.Range("E35").Value = "  +3.24%  "
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").Value = "'2.439"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "'0.5859"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").Value = "'0.01642"
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").Value = "'5.841"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.040.59"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "'0.8398"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "'102.47"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "1.863.84"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  +7.70%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'59.61"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "'8.203"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").Value = "'0.4428"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").Value = "'0.9991"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'0.05277"
$ws.Range("E51").Value = "  -0.77%  "
